$wb = $excel.ActiveWorkbook

# Rename the existing sheet and add the new "v2.5" demo sheet after it
$notes = $wb.ActiveSheet
$notes.Name = "notes"

$demo = $wb.Worksheets.Add($null, $notes)
$demo.Name = "v2.5"

# Populate header row
$demo.Range("A1").Value = "Article"
$demo.Range("B1").Value = "Figure 1"
$demo.Range("C1").Value = "Figure 2"
$demo.Range("D1").Value = "Figure 3"
$demo.Range("E1").Value = "Figure 4"
$demo.Range("F1").Value = "Figure 5"
$demo.Range("G1").Value = "Figure 6"
$demo.Range("H1").Value = "Figure 7"
$demo.Range("I1").Value = "Figure 8"
$demo.Range("J1").Value = "Figure 9"
$demo.Range("K1").Value = "Figure 10"

# Populate data row
$demo.Range("C2").Value = "demo3D27"
$demo.Range("A2").Value = "Johnson et a. (2105)"

# Column width for column A (target stored width 18.33203125 chars)
$demo.Columns.Item(1).ColumnWidth = 17.617745535714285

# Match the page margins used on the "notes" sheet (0.75in/1in/0.5in)
$demo.PageSetup.LeftMargin = 54
$demo.PageSetup.RightMargin = 54
$demo.PageSetup.TopMargin = 72
$demo.PageSetup.BottomMargin = 72
$demo.PageSetup.HeaderMargin = 36
$demo.PageSetup.FooterMargin = 36

$demo.Range("A3").Select() | Out-Null
$demo.Activate() | Out-Null
